$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like a number (e.g. "234.71") or
# like a dotted/grouped figure (e.g. "92.426.17"). The source workbook stores
# all of these as plain text (inlineStr), so force text formatting before
# writing the value, then reset the cell style back to Normal so no stray
# number-format style lingers on the cell.

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '92.426.17'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.112.99'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.42%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '234.71'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -2.73%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '613.62'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("E7").Value = '  -2.61%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.389'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.802'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +7.34%  '

$ws.Range("B11").Value = 'LidoStakedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '3.109.62'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.198'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000244'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -3.41%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '92.197.64'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '33.85'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -3.43%  '

$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '5.42'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -3.14%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '3.690.68'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '3.070.64'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -2.09%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '3.77'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '14.54'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -2.55%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '5.83'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.41%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.0000204'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +1.33%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '9.23'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '439.22'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -3.58%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '5.58'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -5.51%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '85.31'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -3.83%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '11.48'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.68%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '3.267.46'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  +7.51%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.229'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.42%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '0.120'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -20.74%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '1.04'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +45.23%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '9.22'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -1.32%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '7.98'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +7.28%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.157'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -9.79%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '25.87'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -1.89%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '3.90'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -0.63%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '1.90'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -2.91%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '23.87'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +7.83%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '1.28'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -3.00%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '466.02'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -5.11%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.432'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.33%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '3.27'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -3.78%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '159.68'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +2.25%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.684'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -3.10%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.84'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -4.18%  '

$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.32'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -2.49%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.0326'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -0.18%  '

$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '43.84'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -0.47%  '
